$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: nrworkshoop 2 -> 1 (Data unchanged)
$ws.Range("H2").Value = 1

# Row 6: Data "  8/8/2022" -> " 8/11/2022", nrworkshoop 1 -> 2
$ws.Range("G6").Value = " 8/11/2022"
$ws.Range("H6").Value = 2

# Row 8: Data "  9/8/2022" -> " 11/9/2022", nrworkshoop 1 -> 2
$ws.Range("G8").Value = " 11/9/2022"
$ws.Range("H8").Value = 2

# Row 13: nrworkshoop 2 -> 1 (Data unchanged)
$ws.Range("H13").Value = 1

# Row 17: Data "  9/8/2022" -> " 11/9/2022", nrworkshoop 1 -> 2
$ws.Range("G17").Value = " 11/9/2022"
$ws.Range("H17").Value = 2

# Row 21: nrworkshoop 2 -> 1 (Data unchanged)
$ws.Range("H21").Value = 1

# Row 23: Data "  9/8/2022" -> " 11/9/2022", nrworkshoop 1 -> 2
$ws.Range("G23").Value = " 11/9/2022"
$ws.Range("H23").Value = 2

# Row 28: nrworkshoop 2 -> 1 (Data unchanged)
$ws.Range("H28").Value = 1
